# Loan RBI, Variable Instalments
# - Insert a new (blank) column before the old "Late" column (col N) on the
#   "Repayment schedule" sheet, shifting Late / heading / Outstanding one
#   column to the right (N->O, O->P, P->Q), and give the new column a width
#   of 11.
# - Make "Repayment schedule" the active sheet/tab, with selection on R12.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at position N (14th column); this shifts the
# existing N/O/P columns (Late, heading, Outstanding) one place to the right.
[void]$ws.Columns.Item(14).Insert()

# Set the width of the newly inserted column N.
$ws.Columns.Item(14).ColumnWidth = 10.1

# Make "Repayment schedule" the active sheet/tab and select cell R12.
$ws.Activate()
[void]$ws.Range("R12").Select()
